$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.269101333333333
$ws.Range("H2").Value = 6.807304
$ws.Range("I2").Value = 0.02891211995713196
$ws.Range("J2").Value = 0.02891211995713196
$ws.Range("M2").Value = 27.85106533333333
$ws.Range("N2").Value = 83.553196
$ws.Range("O2").Value = 0.1861900221007236
$ws.Range("P2").Value = 0.1861900221007236
$ws.Range("Q2").Value = 63.19688948262045
$ws.Range("R2").Value = 568.7720053435841
$ws.Range("S2").Value = 0.005383148253797171
$ws.Range("T2").Value = 0.005383148253797171
$ws.Range("G3").Value = 2.269101333333333
$ws.Range("H3").Value = 6.807304
$ws.Range("I3").Value = 0.02891211995713196
$ws.Range("J3").Value = 0.02891211995713196
$ws.Range("O3").Value = 0.4727109026912454
$ws.Range("P3").Value = 0.4727109026912454
$ws.Range("Q3").Value = 160.4482256221413
$ws.Range("R3").Value = 1444.034030599272
$ws.Range("S3").Value = 0.01366707432365342
$ws.Range("T3").Value = 0.01366707432365342
$ws.Range("G4").Value = 2.269101333333333
$ws.Range("H4").Value = 6.807304
$ws.Range("I4").Value = 0.02891211995713196
$ws.Range("J4").Value = 0.02891211995713196
$ws.Range("M4").Value = 15.018964
$ws.Range("N4").Value = 45.056892
$ws.Range("O4").Value = 0.1004048213460311
$ws.Range("P4").Value = 0.1004048213460311
$ws.Range("Q4").Value = 34.07955123768533
$ws.Range("R4").Value = 306.715961139168
$ws.Range("S4").Value = 0.002902916239030853
$ws.Range("T4").Value = 0.002902916239030853
$ws.Range("G5").Value = 2.269101333333333
$ws.Range("H5").Value = 6.807304
$ws.Range("I5").Value = 0.02891211995713196
$ws.Range("J5").Value = 0.02891211995713196
$ws.Range("M5").Value = 36.00403133333333
$ws.Range("N5").Value = 108.012094
$ws.Range("O5").Value = 0.2406942538619999
$ws.Range("P5").Value = 0.2406942538619999
$ws.Range("Q5").Value = 81.69679550384177
$ws.Range("R5").Value = 735.271159534576
$ws.Range("S5").Value = 0.006958981140650513
$ws.Range("T5").Value = 0.006958981140650513
$ws.Range("I6").Value = 0.7238963226334669
$ws.Range("J6").Value = 0.7238963226334669
$ws.Range("M6").Value = 27.85106533333333
$ws.Range("N6").Value = 83.553196
$ws.Range("O6").Value = 0.1861900221007236
$ws.Range("P6").Value = 0.1861900221007236
$ws.Range("Q6").Value = 1582.312053428569
$ws.Range("R6").Value = 14240.80848085712
$ws.Range("S6").Value = 0.1347822723097578
$ws.Range("T6").Value = 0.1347822723097578
$ws.Range("I7").Value = 0.7238963226334669
$ws.Range("J7").Value = 0.7238963226334669
$ws.Range("O7").Value = 0.4727109026912454
$ws.Range("P7").Value = 0.4727109026912454
$ws.Range("S7").Value = 0.3421936841269392
$ws.Range("T7").Value = 0.3421936841269392
$ws.Range("I8").Value = 0.7238963226334669
$ws.Range("J8").Value = 0.7238963226334669
$ws.Range("M8").Value = 15.018964
$ws.Range("N8").Value = 45.056892
$ws.Range("O8").Value = 0.1004048213460311
$ws.Range("P8").Value = 0.1004048213460311
$ws.Range("Q8").Value = 853.2775131860813
$ws.Range("R8").Value = 7679.497618674732
$ws.Range("S8").Value = 0.07268268094706211
$ws.Range("T8").Value = 0.07268268094706211
$ws.Range("I9").Value = 0.7238963226334669
$ws.Range("J9").Value = 0.7238963226334669
$ws.Range("M9").Value = 36.00403133333333
$ws.Range("N9").Value = 108.012094
$ws.Range("O9").Value = 0.2406942538619999
$ws.Range("P9").Value = 0.2406942538619999
$ws.Range("Q9").Value = 2045.509285512664
$ws.Range("R9").Value = 18409.58356961397
$ws.Range("S9").Value = 0.1742376852497079
$ws.Range("T9").Value = 0.1742376852497079
$ws.Range("G10").Value = 18.57257166666666
$ws.Range("H10").Value = 55.717715
$ws.Range("I10").Value = 0.2366454120188096
$ws.Range("J10").Value = 0.2366454120188096
$ws.Range("M10").Value = 27.85106533333333
$ws.Range("N10").Value = 83.553196
$ws.Range("O10").Value = 0.1861900221007236
$ws.Range("P10").Value = 0.1861900221007236
$ws.Range("Q10").Value = 517.2659068963488
$ws.Range("R10").Value = 4655.39316206714
$ws.Range("S10").Value = 0.044061014493817
$ws.Range("T10").Value = 0.044061014493817
$ws.Range("G11").Value = 18.57257166666666
$ws.Range("H11").Value = 55.717715
$ws.Range("I11").Value = 0.2366454120188096
$ws.Range("J11").Value = 0.2366454120188096
$ws.Range("O11").Value = 0.4727109026912454
$ws.Range("P11").Value = 0.4727109026912454
$ws.Range("Q11").Value = 1313.267118299721
$ws.Range("R11").Value = 11819.40406469749
$ws.Range("S11").Value = 0.1118648663331532
$ws.Range("T11").Value = 0.1118648663331532
$ws.Range("G12").Value = 18.57257166666666
$ws.Range("H12").Value = 55.717715
$ws.Range("I12").Value = 0.2366454120188096
$ws.Range("J12").Value = 0.2366454120188096
$ws.Range("M12").Value = 15.018964
$ws.Range("N12").Value = 45.056892
$ws.Range("O12").Value = 0.1004048213460311
$ws.Range("P12").Value = 0.1004048213460311
$ws.Range("Q12").Value = 278.9407852490866
$ws.Range("R12").Value = 2510.46706724178
$ws.Range("S12").Value = 0.02376034031610649
$ws.Range("T12").Value = 0.02376034031610649
$ws.Range("G13").Value = 18.57257166666666
$ws.Range("H13").Value = 55.717715
$ws.Range("I13").Value = 0.2366454120188096
$ws.Range("J13").Value = 0.2366454120188096
$ws.Range("M13").Value = 36.00403133333333
$ws.Range("N13").Value = 108.012094
$ws.Range("O13").Value = 0.2406942538619999
$ws.Range("P13").Value = 0.2406942538619999
$ws.Range("Q13").Value = 668.6874522272454
$ws.Range("R13").Value = 6018.187070045209
$ws.Range("S13").Value = 0.05695919087573292
$ws.Range("T13").Value = 0.05695919087573292
$ws.Range("G14").Value = 0.8276899999999999
$ws.Range("H14").Value = 2.48307
$ws.Range("I14").Value = 0.01054614539059158
$ws.Range("J14").Value = 0.01054614539059158
$ws.Range("M14").Value = 27.85106533333333
$ws.Range("N14").Value = 83.553196
$ws.Range("O14").Value = 0.1861900221007236
$ws.Range("P14").Value = 0.1861900221007236
$ws.Range("Q14").Value = 23.05204826574667
$ws.Range("R14").Value = 207.46843439172
$ws.Range("S14").Value = 0.001963587043351691
$ws.Range("T14").Value = 0.001963587043351691
$ws.Range("G15").Value = 0.8276899999999999
$ws.Range("H15").Value = 2.48307
$ws.Range("I15").Value = 0.01054614539059158
$ws.Range("J15").Value = 0.01054614539059158
$ws.Range("O15").Value = 0.4727109026912454
$ws.Range("P15").Value = 0.4727109026912454
$ws.Range("Q15").Value = 58.52598555839
$ws.Range("R15").Value = 526.7338700255099
$ws.Range("S15").Value = 0.004985277907499664
$ws.Range("T15").Value = 0.004985277907499663
$ws.Range("G16").Value = 0.8276899999999999
$ws.Range("H16").Value = 2.48307
$ws.Range("I16").Value = 0.01054614539059158
$ws.Range("J16").Value = 0.01054614539059158
$ws.Range("M16").Value = 15.018964
$ws.Range("N16").Value = 45.056892
$ws.Range("O16").Value = 0.1004048213460311
$ws.Range("P16").Value = 0.1004048213460311
$ws.Range("Q16").Value = 12.43104631316
$ws.Range("R16").Value = 111.87941681844
$ws.Range("S16").Value = 0.001058883843831617
$ws.Range("T16").Value = 0.001058883843831617
$ws.Range("G17").Value = 0.8276899999999999
$ws.Range("H17").Value = 2.48307
$ws.Range("I17").Value = 0.01054614539059158
$ws.Range("J17").Value = 0.01054614539059158
$ws.Range("M17").Value = 36.00403133333333
$ws.Range("N17").Value = 108.012094
$ws.Range("O17").Value = 0.2406942538619999
$ws.Range("P17").Value = 0.2406942538619999
$ws.Range("Q17").Value = 29.80017669428666
$ws.Range("R17").Value = 268.2015902485799
$ws.Range("S17").Value = 0.00253839659590861
$ws.Range("T17").Value = 0.00253839659590861
